$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 69.666664
$ws.Range("I9").Value = 30
$ws.Range("J9").Value = 77.59999999999999
$ws.Range("K9").Value = 30
$ws.Range("L9").Value = 77.59999999999999
$ws.Range("M9").Value = 139
$ws.Range("N9").Value = -415.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3000
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1598.8
$ws.Range("I86").Value = 1498.5
$ws.Range("K86").Value = 1498.5
$ws.Range("M86").Value = -375.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1598.8
$ws.Range("I89").Value = 1498.5
$ws.Range("K89").Value = 7492.5
$ws.Range("M89").Value = -1876.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 481.16666
$ws.Range("I107").Value = 409
$ws.Range("K107").Value = 409
$ws.Range("M107").Value = 1511

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 856.1111
$ws.Range("J129").Value = 880.4103
$ws.Range("L129").Value = 2641.2309
$ws.Range("N129").Value = -12641.2309

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 31754.94
$ws.Range("I137").Value = 1354.2307
$ws.Range("J137").Value = 144671.86
$ws.Range("K137").Value = 4062.6921
$ws.Range("L137").Value = 434015.58
$ws.Range("M137").Value = -1512.6921
$ws.Range("N137").Value = -439115.58

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2873.0466
$ws.Range("J138").Value = 3571.4
$ws.Range("L138").Value = 10714.2
$ws.Range("N138").Value = -20994.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1533.7222
$ws.Range("I45").Value = 975
$ws.Range("J45").Value = 1813.0834
$ws.Range("K45").Value = 975
$ws.Range("L45").Value = 1813.0834
$ws.Range("M45").Value = -598
$ws.Range("N45").Value = -2567.0834

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1321.1875
$ws.Range("I99").Value = 1103.6364
$ws.Range("K99").Value = 1103.6364
$ws.Range("M99").Value = 394.3635999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 232.85715
$ws.Range("I7").Value = 105
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 105
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = 8
$ws.Range("N7").Value = -1226

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 850
$ws.Range("J22").Value = 1300
$ws.Range("L22").Value = 1300
$ws.Range("N22").Value = -2000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2176228
$ws.Range("I58").Value = 3624870
$ws.Range("K58").Value = 3624870
$ws.Range("M58").Value = -3624667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2502478.2
$ws.Range("I99").Value = 9999999
$ws.Range("J99").Value = 3304.6667
$ws.Range("K99").Value = 9999999
$ws.Range("L99").Value = 3304.6667
$ws.Range("M99").Value = -9998501
$ws.Range("N99").Value = -6300.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2502478.2
$ws.Range("I126").Value = 9999999
$ws.Range("J126").Value = 3304.6667
$ws.Range("K126").Value = 29999997
$ws.Range("L126").Value = 9914.000100000001
$ws.Range("M126").Value = -29997527
$ws.Range("N126").Value = -14854.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1228.8928
$ws.Range("I134").Value = 1246.5
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 3739.5
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -1204.5
$ws.Range("N134").Value = -8070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2176228
$ws.Range("I136").Value = 3624870
$ws.Range("K136").Value = 10874610
$ws.Range("M136").Value = -10872060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 12162.223
$ws.Range("I87").Value = 6160
$ws.Range("J87").Value = 24166.666
$ws.Range("K87").Value = 18480
$ws.Range("L87").Value = 72499.99800000001
$ws.Range("M87").Value = -17232
$ws.Range("N87").Value = -74995.99800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 12162.223
$ws.Range("I90").Value = 6160
$ws.Range("J90").Value = 24166.666
$ws.Range("K90").Value = 55440
$ws.Range("L90").Value = 217499.994
$ws.Range("M90").Value = -49200
$ws.Range("N90").Value = -229979.994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3742
$ws.Range("I137").Value = 1936.6666
$ws.Range("J137").Value = 4128.857
$ws.Range("K137").Value = 5809.9998
$ws.Range("L137").Value = 12386.571
$ws.Range("M137").Value = -709.9997999999996
$ws.Range("N137").Value = -22586.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3064.5386
$ws.Range("I138").Value = 2434.1428
$ws.Range("J138").Value = 3800
$ws.Range("K138").Value = 7302.428400000001
$ws.Range("L138").Value = 11400
$ws.Range("M138").Value = -2162.428400000001
$ws.Range("N138").Value = -21680

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 126.26667
$ws.Range("I2").Value = 159.85715
$ws.Range("J2").Value = 96.875
$ws.Range("K2").Value = 159.85715
$ws.Range("L2").Value = 96.875
$ws.Range("M2").Value = -46.85714999999999
$ws.Range("N2").Value = -322.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3951.5715
$ws.Range("I80").Value = 3916
$ws.Range("K80").Value = 3916
$ws.Range("M80").Value = -2918

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3951.5715
$ws.Range("I83").Value = 3916
$ws.Range("K83").Value = 19580
$ws.Range("M83").Value = -14588

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2502.96
$ws.Range("I102").Value = 2503.7
$ws.Range("K102").Value = 2503.7
$ws.Range("M102").Value = -881.6999999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2564.7778
$ws.Range("I46").Value = 1433.3334
$ws.Range("J46").Value = 3130.5
$ws.Range("K46").Value = 1433.3334
$ws.Range("L46").Value = 3130.5
$ws.Range("M46").Value = -1245.3334
$ws.Range("N46").Value = -3506.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1742.7894
$ws.Range("I132").Value = 1562.8667
$ws.Range("J132").Value = 1860.1305
$ws.Range("K132").Value = 4688.6001
$ws.Range("L132").Value = 5580.3915
$ws.Range("M132").Value = -2158.6001
$ws.Range("N132").Value = -10640.3915

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2276.4138
$ws.Range("J136").Value = 4612.375
$ws.Range("L136").Value = 13837.125
$ws.Range("N136").Value = -18937.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1572.125
$ws.Range("I81").Value = 1572.125
$ws.Range("K81").Value = 3144.25
$ws.Range("M81").Value = -2083.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1572.125
$ws.Range("I84").Value = 1572.125
$ws.Range("K84").Value = 15721.25
$ws.Range("M84").Value = -10417.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 549.9286
$ws.Range("I100").Value = 380.9
$ws.Range("K100").Value = 761.8
$ws.Range("M100").Value = -220.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 60666
$ws.Range("J108").Value = 60666
$ws.Range("L108").Value = 60666
$ws.Range("N108").Value = -68346

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 28695
$ws.Range("J119").Value = 28695
$ws.Range("L119").Value = 28695
$ws.Range("N119").Value = -38371

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 12374.917
$ws.Range("I126").Value = 14624.875
$ws.Range("K126").Value = 43874.625
$ws.Range("M126").Value = -41404.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1168.6333
$ws.Range("I132").Value = 925.38464
$ws.Range("K132").Value = 2776.15392
$ws.Range("M132").Value = -246.1539199999997
